# Add a new "batch_no" column in front of the existing candidate table and
# append a second batch (3334) that repeats the first two candidate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts cand_no/last_name/
# first_name/phone_no from A:D to B:E.
$ws.Columns.Item(1).Insert()

# Header for the newly inserted column.
$ws.Range("A1").Value = "batch_no"

# Existing 4 candidate rows (now rows 2-5) belong to batch 3333.
$ws.Range("A2:A5").Value = 3333

# Append two more rows for batch 3334, duplicating the data of the first
# two candidates (rows 2 and 3).
$srcRows = @(2, 3)
$destRow = 6
foreach ($srcRow in $srcRows) {
    $ws.Range("A$destRow").Value = 3334
    $ws.Range("B$destRow").Value = $ws.Range("B$srcRow").Value2
    $ws.Range("C$destRow").Value = $ws.Range("C$srcRow").Value2
    $ws.Range("D$destRow").Value = $ws.Range("D$srcRow").Value2
    $ws.Range("E$destRow").Value = $ws.Range("E$srcRow").Value2
    $ws.Range("E$destRow").NumberFormat = $ws.Range("E$srcRow").NumberFormat
    $destRow++
}
